$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.847426542446107
$ws.Range("D2").Value = 7.765152263472157
$ws.Range("E2").Value = 11.62671632518023
$ws.Range("F2").Value = 38.86257489639515
$ws.Range("G2").Value = 48.7376925318184
$ws.Range("H2").Value = 18.14970985497524
$ws.Range("I2").Value = 27.90380081841142
$ws.Range("J2").Value = 9.205165647024531
$ws.Range("M2").Value = 26.52282352813911
$ws.Range("N2").Value = 17.22607332303565
$ws.Range("C3").Value = 3.865299073906837
$ws.Range("D3").Value = 7.734552171093203
$ws.Range("E3").Value = 11.65696302439233
$ws.Range("F3").Value = 38.52773991519107
$ws.Range("G3").Value = 47.75889101831749
$ws.Range("H3").Value = 18.07413768497346
$ws.Range("I3").Value = 27.62450732682788
$ws.Range("J3").Value = 9.246327743049051
$ws.Range("M3").Value = 25.76207592193492
$ws.Range("N3").Value = 16.99444348550463
$ws.Range("C4").Value = 3.876701577950428
$ws.Range("D4").Value = 7.715548064743717
$ws.Range("E4").Value = 11.67675194266841
$ws.Range("F4").Value = 38.33549911568038
$ws.Range("G4").Value = 47.16878193939282
$ws.Range("H4").Value = 18.03334448388105
$ws.Range("I4").Value = 27.46193935743205
$ws.Range("J4").Value = 9.272768733370652
$ws.Range("M4").Value = 25.28563865001515
$ws.Range("N4").Value = 16.85276405349676
$ws.Range("C5").Value = 3.8814565998973
$ws.Range("D5").Value = 7.707750771790247
$ws.Range("E5").Value = 11.68512305768952
$ws.Range("F5").Value = 38.26058302649887
$ws.Range("G5").Value = 46.931412010912
$ws.Range("H5").Value = 18.01813814095561
$ws.Range("I5").Value = 27.39799860744797
$ws.Range("J5").Value = 9.283838294599605
$ws.Range("M5").Value = 25.08941096158649
$ws.Range("N5").Value = 16.79522959372392
$ws.Range("C6").Value = 3.882252733083782
$ws.Range("D6").Value = 7.706452873784951
$ws.Range("E6").Value = 11.68653164050832
$ws.Range("F6").Value = 38.24835183691386
$ws.Range("G6").Value = 46.89219510602556
$ws.Range("H6").Value = 18.01569889624313
$ws.Range("I6").Value = 27.38752239544769
$ws.Range("J6").Value = 9.285694215207968
$ws.Range("M6").Value = 25.05671044475903
$ws.Range("N6").Value = 16.78569007922955
$ws.Range("C7").Value = 3.876765266085967
$ws.Range("D7").Value = 7.715443120573864
$ws.Range("E7").Value = 11.67686359438621
$ws.Range("F7").Value = 38.33447482670859
$ws.Range("G7").Value = 47.16556763173901
$ws.Range("H7").Value = 18.03313365994642
$ws.Range("I7").Value = 27.4610676084716
$ws.Range("J7").Value = 9.272916826855578
$ws.Range("M7").Value = 25.28300028622513
$ws.Range("N7").Value = 16.85198722390181
$ws.Range("C8").Value = 3.85350038473509
$ws.Range("D8").Value = 7.7546468235126
$ws.Range("E8").Value = 11.63689325325453
$ws.Range("F8").Value = 38.74439569356271
$ws.Range("G8").Value = 48.39816583971859
$ws.Range("H8").Value = 18.12249334538638
$ws.Range("I8").Value = 27.80568941791277
$ws.Range("J8").Value = 9.219116741396929
$ws.Range("M8").Value = 26.26261252018271
$ws.Range("N8").Value = 17.14613354065639
$ws.Range("C9").Value = 3.81125122707112
$ws.Range("D9").Value = 7.829762322439189
$ws.Range("E9").Value = 11.56813135813041
$ws.Range("F9").Value = 39.65085512645666
$ws.Range("G9").Value = 50.88536236269148
$ws.Range("H9").Value = 18.34182554825907
$ws.Range("I9").Value = 28.54918068861758
$ws.Range("J9").Value = 9.122825203310796
$ws.Range("M9").Value = 28.0985840006099
$ws.Range("N9").Value = 17.7244371696016
$ws.Range("C10").Value = 3.782226778772111
$ws.Range("D10").Value = 7.883806399295958
$ws.Range("E10").Value = 11.52342228563093
$ws.Range("F10").Value = 40.37446068086445
$ws.Range("G10").Value = 52.73395321086316
$ws.Range("H10").Value = 18.52917757689767
$ws.Range("I10").Value = 29.13224124394181
$ws.Range("J10").Value = 9.05762114929948
$ws.Range("M10").Value = 29.38230950124097
$ws.Range("N10").Value = 18.14650806683618
$ws.Range("C11").Value = 3.769451782482995
$ws.Range("D11").Value = 7.908129411665392
$ws.Range("E11").Value = 11.5043330132581
$ws.Range("F11").Value = 40.71499308633344
$ws.Range("G11").Value = 53.5751685926453
$ws.Range("H11").Value = 18.61991492488005
$ws.Range("I11").Value = 29.4044655778783
$ws.Range("J11").Value = 9.029145541787534
$ws.Range("M11").Value = 29.94973545531434
$ws.Range("N11").Value = 18.33714821247215
$ws.Range("C12").Value = 3.764675142039299
$ws.Range("D12").Value = 7.917300886550517
$ws.Range("E12").Value = 11.49728311301275
$ws.Range("F12").Value = 40.84547289195646
$ws.Range("G12").Value = 53.8934099665806
$ws.Range("H12").Value = 18.65504865922889
$ws.Range("I12").Value = 29.50846567063181
$ws.Range("J12").Value = 9.018531913501906
$ws.Range("M12").Value = 30.16204878902438
$ws.Range("N12").Value = 18.40908163883411
$ws.Range("C13").Value = 3.76570117521267
$ws.Range("D13").Value = 7.915327408649891
$ws.Range("E13").Value = 11.49879349648504
$ws.Range("F13").Value = 40.81730551404226
$ws.Range("G13").Value = 53.82489009671478
$ws.Range("H13").Value = 18.64744789851304
$ws.Range("I13").Value = 29.4860281585482
$ws.Range("J13").Value = 9.020810228475607
$ws.Range("M13").Value = 30.11643976601178
$ws.Range("N13").Value = 18.39360192053931
$ws.Range("C14").Value = 3.769057587002664
$ws.Range("D14").Value = 7.908884746879753
$ws.Range("E14").Value = 11.50374943519076
$ws.Range("F14").Value = 40.72569774242022
$ws.Range("G14").Value = 53.60135862686697
$ws.Range("H14").Value = 18.62279001081592
$ws.Range("I14").Value = 29.41300391785655
$ws.Range("J14").Value = 9.02826896206612
$ws.Range("M14").Value = 29.96725483992995
$ws.Range("N14").Value = 18.34307173396384
$ws.Range("C15").Value = 3.771121408025921
$ws.Range("D15").Value = 7.904933292195881
$ws.Range("E15").Value = 11.50680835028351
$ws.Range("F15").Value = 40.66978106970275
$ws.Range("G15").Value = 53.46438888080574
$ws.Range("H15").Value = 18.6077864829615
$ws.Range("I15").Value = 29.36839087939
$ws.Range("J15").Value = 9.03285968903916
$ws.Range("M15").Value = 29.87553655921888
$ws.Range("N15").Value = 18.31208515198397
$ws.Range("C16").Value = 3.783070221332039
$ws.Range("D16").Value = 7.882211438394532
$ws.Range("E16").Value = 11.5246948790277
$ws.Range("F16").Value = 40.35242570115866
$ws.Range("G16").Value = 52.67895325069766
$ws.Range("H16").Value = 18.52335697924239
$ws.Range("I16").Value = 29.11458367338309
$ws.Range("J16").Value = 9.059505870679772
$ws.Range("M16").Value = 29.34487818282042
$ws.Range("N16").Value = 18.13401652821288
$ws.Range("C17").Value = 3.790509703253218
$ws.Range("D17").Value = 7.868204200484477
$ws.Range("E17").Value = 11.53598702334394
$ws.Range("F17").Value = 40.16057381557708
$ws.Range("G17").Value = 52.19693248455497
$ws.Range("H17").Value = 18.4729605754283
$ws.Range("I17").Value = 28.96060712512386
$ws.Range("J17").Value = 9.076155442138194
$ws.Range("M17").Value = 29.01496507623512
$ws.Range("N17").Value = 18.0243835574159
$ws.Range("C18").Value = 3.794829054448984
$ws.Range("D18").Value = 7.860123043624467
$ws.Range("E18").Value = 11.54259958769612
$ws.Range("F18").Value = 40.05130159451268
$ws.Range("G18").Value = 51.91973334395085
$ws.Range("H18").Value = 18.44449376281384
$ws.Range("I18").Value = 28.87270706201546
$ws.Range("J18").Value = 9.085843536274671
$ws.Range("M18").Value = 28.82365979866304
$ws.Range("N18").Value = 17.96120095217229
$ws.Range("C19").Value = 3.796298463050982
$ws.Range("D19").Value = 7.85738272818694
$ws.Range("E19").Value = 11.54485871516476
$ws.Range("F19").Value = 40.01449186234611
$ws.Range("G19").Value = 51.82589736957328
$ws.Range("H19").Value = 18.43494522317988
$ws.Range("I19").Value = 28.84306223311318
$ws.Range("J19").Value = 9.089142976274717
$ws.Range("M19").Value = 28.75862700736965
$ws.Range("N19").Value = 17.93978892834078
$ws.Range("C20").Value = 3.789713585329538
$ws.Range("D20").Value = 7.869697847703719
$ws.Range("E20").Value = 11.53477278791002
$ws.Range("F20").Value = 40.18088618167891
$ws.Range("G20").Value = 52.24824219326171
$ws.Range("H20").Value = 18.47827168574768
$ws.Range("I20").Value = 28.97693020376522
$ws.Range("J20").Value = 9.074371514069512
$ws.Range("M20").Value = 29.05024645523181
$ws.Range("N20").Value = 18.03606752804177
$ws.Range("C21").Value = 3.768070077704312
$ws.Range("D21").Value = 7.910778186256096
$ws.Range("E21").Value = 11.50228891006818
$ws.Range("F21").Value = 40.7525645686064
$ws.Range("G21").Value = 53.66702631241848
$ws.Range("H21").Value = 18.63001179693616
$ws.Range("I21").Value = 29.43442880860974
$ws.Range("J21").Value = 9.026073560592756
$ws.Range("M21").Value = 30.01114484516879
$ws.Range("N21").Value = 18.357921156082
$ws.Range("C22").Value = 3.754279900551251
$ws.Range("D22").Value = 7.937397846179345
$ws.Range("E22").Value = 11.48210057526592
$ws.Range("F22").Value = 41.13504208322757
$ws.Range("G22").Value = 54.59233526418089
$ws.Range("H22").Value = 18.73368121264166
$ws.Range("I22").Value = 29.73871932058088
$ws.Range("J22").Value = 8.995495330949799
$ws.Range("M22").Value = 30.62416374085361
$ws.Range("N22").Value = 18.56674094020833
$ws.Range("C23").Value = 3.761607693190004
$ws.Range("D23").Value = 7.923211844771369
$ws.Range("E23").Value = 11.49278042225479
$ws.Range("F23").Value = 40.93013243160159
$ws.Range("G23").Value = 54.09876961154671
$ws.Range("H23").Value = 18.67794593670287
$ws.Range("I23").Value = 29.57585988693591
$ws.Range("J23").Value = 9.011725530609093
$ws.Range("M23").Value = 30.29840956726343
$ws.Range("N23").Value = 18.4554498076834
$ws.Range("C24").Value = 3.790073378460048
$ws.Range("D24").Value = 7.869022657550404
$ws.Range("E24").Value = 11.5353213681689
$ws.Range("F24").Value = 40.17169975724747
$ws.Range("G24").Value = 52.22504529327918
$ws.Range("H24").Value = 18.47586895294529
$ws.Range("I24").Value = 28.96954859225825
$ws.Range("J24").Value = 9.0751776664482
$ws.Range("M24").Value = 29.03430081904557
$ws.Range("N24").Value = 18.03078567722316
$ws.Range("C25").Value = 3.822323781549076
$ws.Range("D25").Value = 7.809637318149696
$ws.Range("E25").Value = 11.58570900583613
$ws.Range("F25").Value = 39.39513275584255
$ws.Range("G25").Value = 50.20725389691098
$ws.Range("H25").Value = 18.27783102345731
$ws.Range("I25").Value = 28.34124383927394
$ws.Range("J25").Value = 9.147896189643925
$ws.Range("M25").Value = 27.61246109862799
$ws.Range("N25").Value = 17.56821997421976